$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 93, pushing existing rows 93-98 down to 95-100.
$ws.Rows.Item(93).Insert()
$ws.Rows.Item(94).Insert()

# --- Row 93 (new) ---
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 45265
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100103
$ws.Range("H93").Value = "Frutos de hueso (carozo)"
$ws.Range("I93").Value = 100103003
$ws.Range("J93").Value = "Damasco"
$ws.Range("K93").Value = "Castle Brite"
$ws.Range("L93").Value = "Especial"
$ws.Range("M93").Value = 100
$ws.Range("N93").Value = 30000
$ws.Range("O93").Value = 30000
$ws.Range("P93").Value = 30000
$ws.Range("Q93").Value = "$/caja 15 kilos"
$ws.Range("R93").Value = "Región de O'Higgins"
$ws.Range("S93").Value = 2000
$ws.Range("T93").Value = 15

# --- Row 94 (new) ---
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45265
$ws.Range("D94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100103
$ws.Range("H94").Value = "Frutos de hueso (carozo)"
$ws.Range("I94").Value = 100103003
$ws.Range("J94").Value = "Damasco"
$ws.Range("K94").Value = "Castle Brite"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 100
$ws.Range("N94").Value = 24000
$ws.Range("O94").Value = 24000
$ws.Range("P94").Value = 24000
$ws.Range("Q94").Value = "$/caja 15 kilos"
$ws.Range("R94").Value = "Región de O'Higgins"
$ws.Range("S94").Value = 1600
$ws.Range("T94").Value = 15
